$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5 - new synthetic data entry (Trades_10m)
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "INT"
$ws.Range("C5").Value = "TRD"
$ws.Range("D5").Value = "BAT"
$ws.Range("E5").Value = "COVID"
$ws.Range("F5").Value = "Trades_10m"
$ws.Range("G5").Value = "Partition Copy"
$ws.Range("H5").Value = "26-04-2020 10:21:54"
$ws.Range("I5").Value = "26-04-2020 10:21:54"
$ws.Range("J5").Value = "In Progress"
$ws.Range("K5").Value = "Export in Progress"
$ws.Range("L5").Value = "Import in Progress"

# Row 6 - new synthetic data entry (Transactions)
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "UAT"
$ws.Range("C6").Value = "COVID"
$ws.Range("D6").Value = "INT"
$ws.Range("E6").Value = "TRD"
$ws.Range("F6").Value = "Transactions"
$ws.Range("G6").Value = "Partition Copy"
$ws.Range("H6").Value = "26-04-2020 20:06:47"
$ws.Range("I6").Value = "26-04-2020 20:06:47"
$ws.Range("J6").Value = "In Progress"
$ws.Range("K6").Value = "Export in Progress"
$ws.Range("L6").Value = "Import in Progress"
